$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.762.31'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = '1.602.66'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').Value = "'212.16"
$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('D9').Value = "'0.247"

$ws.Range('D10').Value = "'19.71"
$ws.Range('E10').Value = '  +0.77%  '

$ws.Range('D11').Value = "'0.0847"
$ws.Range('E11').Value = '  +0.87%  '

$ws.Range('D12').Value = '1.827.60'
$ws.Range('E12').Value = '  +0.25%  '

$ws.Range('D13').Value = '1.593.76'
$ws.Range('E13').Value = '  -0.38%  '

$ws.Range('E14').Value = '  +1.16%  '

$ws.Range('D15').Value = "'0.526"
$ws.Range('E15').Value = '  +0.42%  '

$ws.Range('D16').Value = "'65.03"
$ws.Range('E16').Value = '  -0.13%  '

$ws.Range('D17').Value = '0.0₃0739'
$ws.Range('E17').Value = '  +0.12%  '

$ws.Range('D18').Value = "'209.95"
$ws.Range('E18').Value = '  +0.31%  '

$ws.Range('E19').Value = '  +0.31%  '

$ws.Range('D20').Value = "'7.13"
$ws.Range('E20').Value = '  +1.21%  '

$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('E22').Value = '  -4.53%  '

$ws.Range('D23').Value = "'9.07"
$ws.Range('E23').Value = '  +0.93%  '

$ws.Range('D24').Value = "'143.76"
$ws.Range('E24').Value = '  -0.42%  '

$ws.Range('D25').Value = "'1.00"
$ws.Range('E25').Value = '  +0.11%  '

$ws.Range('D26').Value = "'7.09"
$ws.Range('E26').Value = '  -0.26%  '

$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').Value = "'15.39"
$ws.Range('E28').Value = '  +0.63%  '

$ws.Range('D29').Value = "'0.0508"
$ws.Range('E29').Value = '  -1.52%  '

$ws.Range('E30').Value = '  +0.25%  '

$ws.Range('E31').Value = '  +0.60%  '

$ws.Range('E32').Value = '  +0.48%  '

$ws.Range('D33').Value = '1.293.91'
$ws.Range('E33').Value = '  +0.96%  '

$ws.Range('E34').Value = '  +1.07%  '

$ws.Range('E35').Value = '  +0.32%  '

$ws.Range('E36').Value = '  +14.66%  '

$ws.Range('D37').Value = "'0.596"
$ws.Range('E37').Value = '  -3.79%  '

$ws.Range('E38').Value = '  +0.16%  '

$ws.Range('D39').Value = "'0.832"
$ws.Range('E39').Value = '  -0.17%  '

$ws.Range('D40').Value = "'5.44"
$ws.Range('E40').Value = '  -0.77%  '

$ws.Range('D41').Value = "'2.19"
$ws.Range('E41').Value = '  -0.17%  '

$ws.Range('D42').Value = "'0.781"
$ws.Range('E42').Value = '  -0.32%  '

$ws.Range('D43').Value = "'63.13"
$ws.Range('E43').Value = '  -0.50%  '

$ws.Range('D44').Value = '1.739.72'
$ws.Range('E44').Value = '  +0.51%  '

$ws.Range('D45').Value = "'90.42"
$ws.Range('E45').Value = '  -0.64%  '

$ws.Range('D46').Value = "'6.79"
$ws.Range('E46').Value = '  +32.34%  '

$ws.Range('E47').Value = '  -0.89%  '

$ws.Range('E48').Value = '  +0.72%  '

$ws.Range('D49').Value = "'0.0513"
$ws.Range('E49').Value = '  +0.86%  '

$ws.Range('D50').Value = "'7.56"
$ws.Range('E50').Value = '  +2.14%  '

$ws.Range('E51').Value = '  +0.25%  '
